$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 4.824000000000003
$ws.Range("C4").Value = -14.1
$ws.Range("D4").Value = -7.9056

$ws.Range("C5").Value = -14.7268

$ws.Range("B7").Value = 6.265699999999998

$ws.Range("C8").Value = -12.26669999999999

$ws.Range("D9").Value = -6.569600000000003

$ws.Range("B16").Value = 9.684400000000004
$ws.Range("C16").Value = -12.3488

$ws.Range("D18").Value = -8.509199999999995
